$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# Update email addresses (column C) from "25" to "28"
$ws.Range("C2").Value = "ibu28@gmail.com"
$ws.Range("C3").Value = "deepti.kharbanda28@gmail.com"
$ws.Range("C4").Value = "kartika.varma28@gmail.com"
$ws.Range("C5").Value = "kavita.kharbanda28@gmail.com"

# Update phone numbers (column D) - keep them stored as text, like the originals
$ws.Range("D2").Value = "'1000000054"
$ws.Range("D3").Value = "'1000000055"
$ws.Range("D4").Value = "'1000000056"
$ws.Range("D5").Value = "'1000000057"

# Update selection on the active sheet
$ws.Activate()
$ws.Range("F3").Select()
